$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 131, shifting existing rows 131:224 down to 132:225.
$ws.Rows("131:131").Insert()

# Populate the newly inserted row 131 with the new data record.
$ws.Range("A131").Value = 9
$ws.Range("B131").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C131").Value = "Metropolitana"
$ws.Range("D131").Value = 44873
$ws.Range("E131").Value = 13
$ws.Range("F131").Value = "Fruta"
$ws.Range("G131").Value = 100101
$ws.Range("H131").Value = "Berries"
$ws.Range("I131").Value = 100101001
$ws.Range("J131").Value = "Arándano (blue)"
$ws.Range("K131").Value = "Sin especificar"
$ws.Range("L131").Value = "Primera"
$ws.Range("M131").Value = 450
$ws.Range("N131").Value = 7500
$ws.Range("O131").Value = 8000
$ws.Range("P131").Value = 7778
$ws.Range("Q131").Value = "$/bandeja 2 kilos"
$ws.Range("R131").Value = "Provincia de Linares"
$ws.Range("S131").Value = 3889
$ws.Range("T131").Value = 2
